$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("values")

# Update PANORAMA_IP value (row 4, column B) from 192.168.55.7 to 192.168.55.6
$ws.Cells.Item(4, 2).Value = "192.168.55.6"

# Insert a new row at position 8 (pushes STACK and everything below down by one row)
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the TEMPLATE variable
$ws.Cells.Item(8, 1).Value = "TEMPLATE"
$ws.Cells.Item(8, 2).Value = "sample_template"
$ws.Cells.Item(8, 3).Value = "Template name for Panorama"
